$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Cells.Item(62, 8).Value = 1033.3334  # H62: 1090 -> 1033.3334
$ws.Cells.Item(62, 9).Value = 975  # I62: 1050 -> 975
$ws.Cells.Item(62, 10).Value = 1091.6666  # J62: 1116.6666 -> 1091.6666
$ws.Cells.Item(62, 11).Value = 975  # K62: 1050 -> 975
$ws.Cells.Item(62, 12).Value = 1091.6666  # L62: 1116.6666 -> 1091.6666
$ws.Cells.Item(62, 13).Value = -351  # M62: -426 -> -351
$ws.Cells.Item(62, 14).Value = -2339.6666  # N62: -2364.6666 -> -2339.6666
# row 64
$ws.Cells.Item(64, 8).Value = 336000  # H64: 46590.74 -> 336000
$ws.Cells.Item(64, 9).Value = 1000000  # I64: 335333.34 -> 1000000
$ws.Cells.Item(64, 10).Value = 4000  # J64: 3279.35 -> 4000
$ws.Cells.Item(64, 11).Value = 1000000  # K64: 335333.34 -> 1000000
$ws.Cells.Item(64, 12).Value = 4000  # L64: 3279.35 -> 4000
$ws.Cells.Item(64, 13).Value = -999752  # M64: -335085.34 -> -999752
$ws.Cells.Item(64, 14).Value = -4496  # N64: -3775.35 -> -4496
# row 65
$ws.Cells.Item(65, 8).Value = 1033.3334  # H65: 1090 -> 1033.3334
$ws.Cells.Item(65, 9).Value = 975  # I65: 1050 -> 975
$ws.Cells.Item(65, 10).Value = 1091.6666  # J65: 1116.6666 -> 1091.6666
$ws.Cells.Item(65, 11).Value = 4875  # K65: 5250 -> 4875
$ws.Cells.Item(65, 12).Value = 5458.333000000001  # L65: 5583.333000000001 -> 5458.333000000001
$ws.Cells.Item(65, 13).Value = -1755  # M65: -2130 -> -1755
$ws.Cells.Item(65, 14).Value = -11698.333  # N65: -11823.333 -> -11698.333
# row 67
$ws.Cells.Item(67, 8).Value = 336000  # H67: 46590.74 -> 336000
$ws.Cells.Item(67, 9).Value = 1000000  # I67: 335333.34 -> 1000000
$ws.Cells.Item(67, 10).Value = 4000  # J67: 3279.35 -> 4000
$ws.Cells.Item(67, 11).Value = 1000000  # K67: 335333.34 -> 1000000
$ws.Cells.Item(67, 12).Value = 4000  # L67: 3279.35 -> 4000
$ws.Cells.Item(67, 13).Value = -999142  # M67: -334475.34 -> -999142
$ws.Cells.Item(67, 14).Value = -5716  # N67: -4995.35 -> -5716
# row 74
$ws.Cells.Item(74, 8).Value = 5000  # H74: 4150 -> 5000
$ws.Cells.Item(74, 9).Value = 0  # I74: 3300 -> 0
$ws.Cells.Item(74, 11).Value = 0  # K74: 3300 -> 0
$ws.Cells.Item(74, 13).ClearContents()  # M74: -2364 -> (removed)
# row 76
$ws.Cells.Item(76, 8).Value = 4507.846  # H76: 4100.8335 -> 4507.846
$ws.Cells.Item(76, 9).Value = 3966.6667  # I76: 3451.5 -> 3966.6667
$ws.Cells.Item(76, 10).Value = 4670.2  # J76: 4182 -> 4670.2
$ws.Cells.Item(76, 11).Value = 3966.6667  # K76: 3451.5 -> 3966.6667
$ws.Cells.Item(76, 12).Value = 4670.2  # L76: 4182 -> 4670.2
$ws.Cells.Item(76, 13).Value = -3651.6667  # M76: -3136.5 -> -3651.6667
$ws.Cells.Item(76, 14).Value = -5300.2  # N76: -4812 -> -5300.2
# row 77
$ws.Cells.Item(77, 8).Value = 5000  # H77: 4150 -> 5000
$ws.Cells.Item(77, 9).Value = 0  # I77: 3300 -> 0
$ws.Cells.Item(77, 11).Value = 0  # K77: 16500 -> 0
$ws.Cells.Item(77, 13).ClearContents()  # M77: -11820 -> (removed)
# row 79
$ws.Cells.Item(79, 8).Value = 4507.846  # H79: 4100.8335 -> 4507.846
$ws.Cells.Item(79, 9).Value = 3966.6667  # I79: 3451.5 -> 3966.6667
$ws.Cells.Item(79, 10).Value = 4670.2  # J79: 4182 -> 4670.2
$ws.Cells.Item(79, 11).Value = 3966.6667  # K79: 3451.5 -> 3966.6667
$ws.Cells.Item(79, 12).Value = 4670.2  # L79: 4182 -> 4670.2
$ws.Cells.Item(79, 13).Value = -2874.6667  # M79: -2359.5 -> -2874.6667
$ws.Cells.Item(79, 14).Value = -6854.2  # N79: -6366 -> -6854.2
# row 99
$ws.Cells.Item(99, 8).Value = 12798.5  # H99: 11402.777 -> 12798.5
$ws.Cells.Item(99, 9).Value = 14484  # I99: 12703.125 -> 14484
$ws.Cells.Item(99, 11).Value = 43452  # K99: 38109.375 -> 43452
$ws.Cells.Item(99, 13).Value = -41954  # M99: -36611.375 -> -41954
# row 112
$ws.Cells.Item(112, 8).Value = 1125.6842  # H112: 1230.3077 -> 1125.6842
$ws.Cells.Item(112, 10).Value = 1160.4445  # J112: 1291.1666 -> 1160.4445
$ws.Cells.Item(112, 12).Value = 3481.3335  # L112: 3873.4998 -> 3481.3335
$ws.Cells.Item(112, 14).Value = -5697.333500000001  # N112: -6089.4998 -> -5697.333500000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 8
$ws.Cells.Item(8, 8).Value = 49400  # H8: 36266.668 -> 49400
$ws.Cells.Item(8, 10).Value = 49400  # J8: 36266.668 -> 49400
$ws.Cells.Item(8, 12).Value = 49400  # L8: 36266.668 -> 49400
$ws.Cells.Item(8, 14).Value = -49688  # N8: -36554.668 -> -49688
# row 32
$ws.Cells.Item(32, 8).Value = 24497.24  # H32: 20730.686 -> 24497.24
$ws.Cells.Item(32, 9).Value = 4297.2344  # I32: 3594.443 -> 4297.2344
$ws.Cells.Item(32, 10).Value = 142024.55  # J32: 156107 -> 142024.55
$ws.Cells.Item(32, 11).Value = 4297.2344  # K32: 3594.443 -> 4297.2344
$ws.Cells.Item(32, 12).Value = 142024.55  # L32: 156107 -> 142024.55
$ws.Cells.Item(32, 13).Value = -4010.2344  # M32: -3307.443 -> -4010.2344
$ws.Cells.Item(32, 14).Value = -142598.55  # N32: -156681 -> -142598.55
# row 61
$ws.Cells.Item(61, 8).Value = 1938.9584  # H61: 2139.476 -> 1938.9584
$ws.Cells.Item(61, 9).Value = 1408.6923  # I61: 1631.619 -> 1408.6923
$ws.Cells.Item(61, 10).Value = 2565.6365  # J61: 2647.3333 -> 2565.6365
$ws.Cells.Item(61, 11).Value = 1408.6923  # K61: 1631.619 -> 1408.6923
$ws.Cells.Item(61, 12).Value = 2565.6365  # L61: 2647.3333 -> 2565.6365
$ws.Cells.Item(61, 13).Value = -1196.6923  # M61: -1419.619 -> -1196.6923
$ws.Cells.Item(61, 14).Value = -2989.6365  # N61: -3071.3333 -> -2989.6365
# row 74
$ws.Cells.Item(74, 8).Value = 2632850  # H74: 2942607.2 -> 2632850
$ws.Cells.Item(74, 9).Value = 836.14813  # I74: 962.1739 -> 836.14813
$ws.Cells.Item(74, 10).Value = 9093248  # J74: 9093319 -> 9093248
$ws.Cells.Item(74, 11).Value = 836.14813  # K74: 962.1739 -> 836.14813
$ws.Cells.Item(74, 12).Value = 9093248  # L74: 9093319 -> 9093248
$ws.Cells.Item(74, 13).Value = 37.85186999999996  # M74: -88.1739 -> 37.85186999999996
$ws.Cells.Item(74, 14).Value = -9094996  # N74: -9095067 -> -9094996
# row 77
$ws.Cells.Item(77, 8).Value = 2632850  # H77: 2942607.2 -> 2632850
$ws.Cells.Item(77, 9).Value = 836.14813  # I77: 962.1739 -> 836.14813
$ws.Cells.Item(77, 10).Value = 9093248  # J77: 9093319 -> 9093248
$ws.Cells.Item(77, 11).Value = 4180.74065  # K77: 4810.8695 -> 4180.74065
$ws.Cells.Item(77, 12).Value = 45466240  # L77: 45466595 -> 45466240
$ws.Cells.Item(77, 13).Value = 187.2593500000003  # M77: -442.8694999999998 -> 187.2593500000003
$ws.Cells.Item(77, 14).Value = -45474976  # N77: -45475331 -> -45474976
# row 132
$ws.Cells.Item(132, 8).Value = 2054.8635  # H132: 2479.0286 -> 2054.8635
$ws.Cells.Item(132, 9).Value = 1876  # I132: 2289.625 -> 1876
$ws.Cells.Item(132, 11).Value = 5628  # K132: 6868.875 -> 5628
$ws.Cells.Item(132, 13).Value = -3098  # M132: -4338.875 -> -3098
# row 136
$ws.Cells.Item(136, 8).Value = 1938.9584  # H136: 2139.476 -> 1938.9584
$ws.Cells.Item(136, 9).Value = 1408.6923  # I136: 1631.619 -> 1408.6923
$ws.Cells.Item(136, 10).Value = 2565.6365  # J136: 2647.3333 -> 2565.6365
$ws.Cells.Item(136, 11).Value = 4226.0769  # K136: 4894.857 -> 4226.0769
$ws.Cells.Item(136, 12).Value = 7696.9095  # L136: 7941.999899999999 -> 7696.9095
$ws.Cells.Item(136, 13).Value = -1676.0769  # M136: -2344.857 -> -1676.0769
$ws.Cells.Item(136, 14).Value = -12796.9095  # N136: -13041.9999 -> -12796.9095

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Cells.Item(20, 8).Value = 2060.3157  # H20: 41408 -> 2060.3157
$ws.Cells.Item(20, 9).Value = 2133.5625  # I20: 49109.523 -> 2133.5625
$ws.Cells.Item(20, 10).Value = 1669.6666  # J20: 975 -> 1669.6666
$ws.Cells.Item(20, 11).Value = 2133.5625  # K20: 49109.523 -> 2133.5625
$ws.Cells.Item(20, 12).Value = 1669.6666  # L20: 975 -> 1669.6666
$ws.Cells.Item(20, 13).Value = -1886.5625  # M20: -48862.523 -> -1886.5625
$ws.Cells.Item(20, 14).Value = -2163.6666  # N20: -1469 -> -2163.6666

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 112
$ws.Cells.Item(112, 8).Value = 92957  # H112: 144928.58 -> 92957
$ws.Cells.Item(112, 9).Value = 334175.66  # I112: 501000 -> 334175.66
$ws.Cells.Item(112, 11).Value = 1002526.98  # K112: 1503000 -> 1002526.98
$ws.Cells.Item(112, 13).Value = -1001418.98  # M112: -1501892 -> -1001418.98

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 11
$ws.Cells.Item(11, 8).Value = 4376470.5  # H11: 5237000 -> 4376470.5
$ws.Cells.Item(11, 9).Value = 6091833.5  # I11: 6644545.5 -> 6091833.5
$ws.Cells.Item(11, 10).Value = 259599.4  # J11: 75999 -> 259599.4
$ws.Cells.Item(11, 11).Value = 6091833.5  # K11: 6644545.5 -> 6091833.5
$ws.Cells.Item(11, 12).Value = 259599.4  # L11: 75999 -> 259599.4
$ws.Cells.Item(11, 13).Value = -6091694.5  # M11: -6644406.5 -> -6091694.5
$ws.Cells.Item(11, 14).Value = -259877.4  # N11: -76277 -> -259877.4
# row 23
$ws.Cells.Item(23, 8).Value = 17375.125  # H23: 19800 -> 17375.125
$ws.Cells.Item(23, 9).Value = 1  # I23: 0 -> 1
$ws.Cells.Item(23, 10).Value = 23166.5  # J23: 19800 -> 23166.5
$ws.Cells.Item(23, 11).Value = 1  # K23: 0 -> 1
$ws.Cells.Item(23, 12).Value = 23166.5  # L23: 19800 -> 23166.5
$ws.Cells.Item(23, 13).Value = 222  # M23: None -> 222
$ws.Cells.Item(23, 14).Value = -23612.5  # N23: -20246 -> -23612.5
# row 80
$ws.Cells.Item(80, 8).Value = 333345000  # H80: 250251870 -> 333345000
$ws.Cells.Item(80, 9).Value = 333345000  # I80: 250251870 -> 333345000
$ws.Cells.Item(80, 11).Value = 333345000  # K80: 250251870 -> 333345000
$ws.Cells.Item(80, 13).Value = -333344002  # M80: -250250872 -> -333344002
# row 83
$ws.Cells.Item(83, 8).Value = 333345000  # H83: 250251870 -> 333345000
$ws.Cells.Item(83, 9).Value = 333345000  # I83: 250251870 -> 333345000
$ws.Cells.Item(83, 11).Value = 1666725000  # K83: 1251259350 -> 1666725000
$ws.Cells.Item(83, 13).Value = -1666720008  # M83: -1251254358 -> -1666720008

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 123
$ws.Cells.Item(123, 8).Value = 39000  # H123: 0 -> 39000
$ws.Cells.Item(123, 10).Value = 39000  # J123: 0 -> 39000
$ws.Cells.Item(123, 12).Value = 39000  # L123: 0 -> 39000
$ws.Cells.Item(123, 14).Value = -48800  # N123: None -> -48800
# row 132
$ws.Cells.Item(132, 8).Value = 2539.3333  # H132: 2284.8572 -> 2539.3333
$ws.Cells.Item(132, 9).Value = 2608.2888  # I132: 2353.327 -> 2608.2888
$ws.Cells.Item(132, 10).Value = 2332.4666  # J132: 2087.0557 -> 2332.4666
$ws.Cells.Item(132, 11).Value = 7824.866399999999  # K132: 7059.981000000001 -> 7824.866399999999
$ws.Cells.Item(132, 12).Value = 6997.399800000001  # L132: 6261.1671 -> 6997.399800000001
$ws.Cells.Item(132, 13).Value = -5294.866399999999  # M132: -4529.981000000001 -> -5294.866399999999
$ws.Cells.Item(132, 14).Value = -12057.3998  # N132: -11321.1671 -> -12057.3998

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 30
$ws.Cells.Item(30, 8).Value = 25003  # H30: 35005 -> 25003
$ws.Cells.Item(30, 9).Value = 25003  # I30: 60009 -> 25003
$ws.Cells.Item(30, 10).Value = 0  # J30: 10001 -> 0
$ws.Cells.Item(30, 11).Value = 25003  # K30: 60009 -> 25003
$ws.Cells.Item(30, 12).Value = 0  # L30: 10001 -> 0
$ws.Cells.Item(30, 13).Value = -24896  # M30: -59902 -> -24896
$ws.Cells.Item(30, 14).ClearContents()  # N30: -10215 -> (removed)
# row 132
$ws.Cells.Item(132, 8).Value = 2705.025  # H132: 3004.6943 -> 2705.025
$ws.Cells.Item(132, 9).Value = 2532.4614  # I132: 2913.3635 -> 2532.4614
$ws.Cells.Item(132, 10).Value = 3025.5  # J132: 3148.2144 -> 3025.5
$ws.Cells.Item(132, 11).Value = 7597.3842  # K132: 8740.0905 -> 7597.3842
$ws.Cells.Item(132, 12).Value = 9076.5  # L132: 9444.643199999999 -> 9076.5
$ws.Cells.Item(132, 13).Value = -5067.3842  # M132: -6210.0905 -> -5067.3842
$ws.Cells.Item(132, 14).Value = -14136.5  # N132: -14504.6432 -> -14136.5
# row 136
$ws.Cells.Item(136, 8).Value = 745.5333000000001  # H136: 699.4693600000001 -> 745.5333000000001
$ws.Cells.Item(136, 9).Value = 611.5806  # I136: 559.6857 -> 611.5806
$ws.Cells.Item(136, 10).Value = 1042.1428  # J136: 1048.9286 -> 1042.1428
$ws.Cells.Item(136, 11).Value = 1834.7418  # K136: 1679.0571 -> 1834.7418
$ws.Cells.Item(136, 12).Value = 3126.4284  # L136: 3146.7858 -> 3126.4284
$ws.Cells.Item(136, 13).Value = 715.2582  # M136: 870.9429 -> 715.2582
$ws.Cells.Item(136, 14).Value = -8226.428400000001  # N136: -8246.7858 -> -8226.428400000001
